$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had 10 columns (A:J):
#   A dispositivo | B sesion_id | C dia | D mes | E año | F timestamp
#   G ubicacion   | H modelo    | I variable | J valor
#
# New layout has 8 columns (A:H):
#   A patente | B sesion_id | C dia | D mes | E año | F timestamp
#   G variable | H valor
#
# Remove old columns H (modelo) and I (variable); this shifts the old J
# (valor) column left so it becomes the new column H.
$ws.Range("H1:I20").Delete(-4159)  # xlShiftToLeft

# Update header labels first
$ws.Range("A1").Value = "patente"
$ws.Range("G1").Value = "variable"

# Then update the data values
$ws.Range("G2:G20").Value = "°C"
$ws.Range("A2:A20").Value = "MP-01-EXPRESS"

$ws.Range("K20").Select() | Out-Null
